$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main_Input")
$ws.Range("D2").Value = 0.75
